{"js": "// Update the \"three-digit number divided by one-digit number\" practice\n// sheet: replace each division prompt's text with its new value. Pairs are\n// listed in the same left-to-right / top-to-bottom order they appear in the\n// document (and in the diff).\nconst replacements = [\n  [\"390\u00f75=\", \"117\u00f78=\"],\n  [\"328\u00f73=\", \"582\u00f74=\"],\n  [\"766\u00f73=\", \"858\u00f77=\"],\n  [\"266\u00f72=\", \"341\u00f75=\"],\n  [\"280\u00f77=\", \"510\u00f77=\"],\n  [\"143\u00f78=\", \"813\u00f72=\"],\n  [\"462\u00f78=\", \"409\u00f79=\"],\n  [\"320\u00f79=\", \"365\u00f75=\"],\n  [\"594\u00f73=\", \"804\u00f73=\"],\n  [\"362\u00f72=\", \"548\u00f78=\"],\n  [\"400\u00f74=\", \"646\u00f74=\"],\n  [\"850\u00f73=\", \"726\u00f79=\"],\n  [\"361\u00f72=\", \"514\u00f74=\"],\n  [\"288\u00f76=\", \"182\u00f73=\"],\n  [\"906\u00f78=\", \"313\u00f75=\"],\n  [\"183\u00f78=\", \"459\u00f72=\"],\n  [\"841\u00f76=\", \"433\u00f78=\"],\n  [\"826\u00f77=\", \"963\u00f78=\"],\n  [\"409\u00f79=\", \"562\u00f72=\"],\n  [\"364\u00f76=\", \"796\u00f79=\"],\n  [\"993\u00f73=\", \"350\u00f77=\"],\n  [\"951\u00f78=\", \"872\u00f77=\"],\n  [\"655\u00f79=\", \"475\u00f76=\"],\n  [\"619\u00f73=\", \"615\u00f72=\"],\n  [\"982\u00f75=\", \"865\u00f75=\"],\n];\n\n// A plain whole-document text search-and-replace is unsafe here: several\n// \"after\" values equal OTHER pairs' \"before\" values (e.g. pair #7 turns\n// \"462\u00f78=\" into \"409\u00f79=\", while pair #19's source is itself \"409\u00f79=\"). If\n// replacements ran as independent body-wide searches, pair #19 could match\n// the text that pair #7 just wrote instead of the original paragraph that\n// held \"409\u00f79=\" before any edits. To keep each pair scoped to the exact\n// paragraph it targets, walk the body's paragraphs once, collect the ones\n// whose (pre-edit) text is one of our known \"before\" values, and pair them\n// up with the replacement list by their left-to-right document order\n// (which matches the order the values appear in the diff/commit).\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\n\nconst beforeSet = new Set(replacements.map(([oldText]) => oldText));\nconst targets = [];\nfor (let i = 0; i < paras.items.length; i++) {\n  if (beforeSet.has(paras.items[i].text)) {\n    targets.push(paras.items[i]);\n  }\n}\n\nif (targets.length !== replacements.length) {\n  throw new Error(\n    `expected ${replacements.length} matching paragraphs, found ${targets.length}`\n  );\n}\n\nfor (let i = 0; i < targets.length; i++) {\n  const [, newText] = replacements[i];\n  targets[i].getRange().insertText(newText, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Update the \"three-digit number divided by one-digit number\" practice\n# sheet: replace each division prompt's text with its new value. Pairs are\n# listed in the same left-to-right / top-to-bottom order they appear in the\n# document (and in the diff).\n$d = $word.ActiveDocument\n\n$oldValues = @(\n    \"390\u00f75=\", \"328\u00f73=\", \"766\u00f73=\", \"266\u00f72=\", \"280\u00f77=\",\n    \"143\u00f78=\", \"462\u00f78=\", \"320\u00f79=\", \"594\u00f73=\", \"362\u00f72=\",\n    \"400\u00f74=\", \"850\u00f73=\", \"361\u00f72=\", \"288\u00f76=\", \"906\u00f78=\",\n    \"183\u00f78=\", \"841\u00f76=\", \"826\u00f77=\", \"409\u00f79=\", \"364\u00f76=\",\n    \"993\u00f73=\", \"951\u00f78=\", \"655\u00f79=\", \"619\u00f73=\", \"982\u00f75=\"\n)\n$newValues = @(\n    \"117\u00f78=\", \"582\u00f74=\", \"858\u00f77=\", \"341\u00f75=\", \"510\u00f77=\",\n    \"813\u00f72=\", \"409\u00f79=\", \"365\u00f75=\", \"804\u00f73=\", \"548\u00f78=\",\n    \"646\u00f74=\", \"726\u00f79=\", \"514\u00f74=\", \"182\u00f73=\", \"313\u00f75=\",\n    \"459\u00f72=\", \"433\u00f78=\", \"963\u00f78=\", \"562\u00f72=\", \"796\u00f79=\",\n    \"350\u00f77=\", \"872\u00f77=\", \"475\u00f76=\", \"615\u00f72=\", \"865\u00f75=\"\n)\n\n# A plain whole-document Find/Replace (Find.Execute ... Replace:=wdReplaceAll)\n# is unsafe here: several \"after\" values equal OTHER pairs' \"before\" values\n# (e.g. pair #7 turns \"462\u00f78=\" into \"409\u00f79=\", while pair #19's source is\n# itself \"409\u00f79=\"). Running independent whole-document replacements risks a\n# later pair matching text an earlier pair just wrote. Instead, walk the\n# document's paragraphs once, collect the ones whose (pre-edit) text is one\n# of our known \"before\" values, and pair them up with the replacement list\n# by their in-document order (which matches the order the values appear in\n# the diff/commit) \u2014 each paragraph is then addressed directly, so no\n# cross-matching between pairs can happen.\n$targets = New-Object System.Collections.ArrayList\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    $t2 = $t.TrimEnd([char]0x07, [char]0x0D)\n    if ($oldValues -contains $t2) {\n        [void]$targets.Add($p)\n    }\n}\n\nif ($targets.Count -ne $oldValues.Length) {\n    throw \"expected $($oldValues.Length) matching paragraphs, found $($targets.Count)\"\n}\n\nfor ($i = 0; $i -lt $targets.Count; $i++) {\n    $targets[$i].Range.Text = $newValues[$i]\n}\n"}
